# Applies the row-level edits to the order matrix sheet for subject 12 / A_block2 / VR
# (5-second black screens replaced by the 'mareo.mp4' motion-sickness clip, and the
# block 1-4 stimulus/instruction ordering columns D-L updated to match the new sequence.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = './instructions_videos/block_2_text.mp4'
$ws.Range("I2").Value = 2

# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '11'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = 'inverse'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '../stimuli/exp_videos/VR/11.mp4'
$ws.Range("I3").Value = 2

# Row 4
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = './instructions_videos/post_stimulus_self_report.mp4'
$ws.Range("I4").Value = 2
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = 'post_stimulus_self_report'

# Row 5
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = './instructions_videos/mareo.mp4'
$ws.Range("I5").Value = 2
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = 'motion_sickness'

# Row 6
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = './instructions_videos/block_2_text_reminder.mp4'
$ws.Range("I6").Value = 2
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = 'audio_instruction'

# Row 7
$ws.Range("D7").Value = 2
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5'
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = 'arousal'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = 'inverse'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '../stimuli/exp_videos/VR/5.mp4'
$ws.Range("I7").Value = 2
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = 'video'

# Row 8
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = './instructions_videos/post_stimulus_self_report.mp4'
$ws.Range("I8").Value = 2
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = 'post_stimulus_self_report'

# Row 9
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = './instructions_videos/mareo.mp4'
$ws.Range("I9").Value = 2
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = 'motion_sickness'

# Row 10
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = './instructions_videos/block_2_text_reminder.mp4'
$ws.Range("I10").Value = 2
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = 'audio_instruction'

# Row 11
$ws.Range("D11").Value = 3
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '10'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = 'arousal'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = 'inverse'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '../stimuli/exp_videos/VR/10.mp4'
$ws.Range("I11").Value = 2
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = 'video'

# Row 12
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = './instructions_videos/post_stimulus_self_report.mp4'
$ws.Range("I12").Value = 2
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = 'post_stimulus_self_report'

# Row 13
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = './instructions_videos/mareo.mp4'
$ws.Range("I13").Value = 2
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value = 'motion_sickness'

# Row 14
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = './instructions_videos/block_2_text_reminder.mp4'
$ws.Range("I14").Value = 2

# Row 15
$ws.Range("D15").Value = 4
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = 'inverse '
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '../stimuli/exp_videos/VR/1.mp4'
$ws.Range("I15").Value = 2

# Row 16
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = './instructions_videos/post_stimulus_self_report.mp4'
$ws.Range("I16").Value = 2
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = 'post_stimulus_self_report'

# Row 17
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = './instructions_videos/mareo.mp4'
$ws.Range("I17").Value = 2
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = 'motion_sickness'

# Row 18
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = './instructions_videos/luminance_instructions_inverse.mp4'
$ws.Range("I18").Value = ""
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = 'luminance_instructions'

# Row 19
$ws.Range("D19").Value = 5
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = 'luminance'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = 'inverse'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '../stimuli/exp_videos/VR/green_intensity_video_7.mp4'
$ws.Range("I19").Value = ""
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = 'luminance'

# Row 20
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = './instructions_videos/confidence_luminance_practice_instructions_text.mp4'
$ws.Range("I20").Value = ""
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L20").Value = 'confidence_luminance_instructions'

# Row 21
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = './instructions_videos/block_3_text.mp4'
$ws.Range("I21").Value = 3
$ws.Range("J21").Value = 3
$ws.Range("L21").NumberFormat = "@"
$ws.Range("L21").Value = 'audio_instruction'

# Row 22
$ws.Range("D22").Value = 6
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '6'
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = 'valence'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = 'direct'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = '../stimuli/exp_videos/VR/6.mp4'
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = 'video'

# Row 23
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = './instructions_videos/post_stimulus_verbal_report.mp4'
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 3
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = 'instruction_post_stimulus_verbal_report'

# Row 24
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = './videos_fixation/countdown_bar.mp4'
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 3
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value = 'verbal_report'

# Row 25
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = './instructions_videos/confidence_verbal_report_text.mp4'
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 3
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = 'confidence_verbal_report'

# Row 26
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = './instructions_videos/mareo.mp4'
$ws.Range("I26").Value = 3
$ws.Range("L26").NumberFormat = "@"
$ws.Range("L26").Value = 'motion_sickness'

# Row 27
$ws.Range("D27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("G27").Value = ""
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = './instructions_videos/block_3_text_reminder.mp4'
$ws.Range("I27").Value = 3
$ws.Range("L27").NumberFormat = "@"
$ws.Range("L27").Value = 'audio_instruction'

# Row 28
$ws.Range("D28").Value = 7
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '13'
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = 'valence'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = 'direct'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '../stimuli/exp_videos/VR/13.mp4'
$ws.Range("I28").Value = 3
$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = 'video'

# Row 29
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = ""
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = './instructions_videos/post_stimulus_verbal_report.mp4'
$ws.Range("I29").Value = 3
$ws.Range("L29").NumberFormat = "@"
$ws.Range("L29").Value = 'instruction_post_stimulus_verbal_report'

# Row 30
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = './videos_fixation/countdown_bar.mp4'
$ws.Range("I30").Value = 3
$ws.Range("L30").NumberFormat = "@"
$ws.Range("L30").Value = 'verbal_report'

# Row 31
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = './instructions_videos/confidence_verbal_report_text.mp4'
$ws.Range("I31").Value = 3
$ws.Range("L31").NumberFormat = "@"
$ws.Range("L31").Value = 'confidence_verbal_report'

# Row 32
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = './instructions_videos/mareo.mp4'
$ws.Range("I32").Value = 3
$ws.Range("L32").NumberFormat = "@"
$ws.Range("L32").Value = 'motion_sickness'

# Row 33
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = ""
$ws.Range("G33").Value = ""
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = './instructions_videos/block_3_text_reminder.mp4'
$ws.Range("I33").Value = 3
$ws.Range("L33").NumberFormat = "@"
$ws.Range("L33").Value = 'audio_instruction'

# Row 34
$ws.Range("D34").Value = 8
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '14'
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = 'valence'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = 'direct'
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '../stimuli/exp_videos/VR/14.mp4'
$ws.Range("I34").Value = 3
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = 'video'

# Row 35
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = './instructions_videos/post_stimulus_verbal_report.mp4'
$ws.Range("I35").Value = 3
$ws.Range("L35").NumberFormat = "@"
$ws.Range("L35").Value = 'instruction_post_stimulus_verbal_report'

# Row 36
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = './videos_fixation/countdown_bar.mp4'
$ws.Range("I36").Value = 3
$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = 'verbal_report'

# Row 37
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = ""
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = ""
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = './instructions_videos/confidence_verbal_report_text.mp4'
$ws.Range("I37").Value = 3
$ws.Range("L37").NumberFormat = "@"
$ws.Range("L37").Value = 'confidence_verbal_report'

# Row 38
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = './instructions_videos/mareo.mp4'
$ws.Range("I38").Value = 3
$ws.Range("L38").NumberFormat = "@"
$ws.Range("L38").Value = 'motion_sickness'

# Row 39
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = './instructions_videos/luminance_instructions_direct.mp4'
$ws.Range("I39").Value = ""
$ws.Range("L39").NumberFormat = "@"
$ws.Range("L39").Value = 'luminance_instructions'

# Row 40
$ws.Range("D40").Value = 9
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = 'luminance'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = 'direct'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '../stimuli/exp_videos/VR/green_intensity_video_9.mp4'
$ws.Range("L40").NumberFormat = "@"
$ws.Range("L40").Value = 'luminance'

# Row 41
$ws.Range("D41").Value = ""
$ws.Range("F41").Value = ""
$ws.Range("G41").Value = ""
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = './instructions_videos/confidence_luminance_practice_instructions_text.mp4'
$ws.Range("L41").NumberFormat = "@"
$ws.Range("L41").Value = 'confidence_luminance_instructions'
